$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Mark the new risk rows (13-15) as Status = "new" ---
$ws.Range("A13").Value = "new"
$ws.Range("A14").Value = "new"
$ws.Range("A15").Value = "new"

# --- 2. Fill in the new row 12 risk entry ---
$ws.Range("A12").Value = "open"

# --- 3. Apply an AutoFilter on the Status column of the risk register table ---
# Only rows whose Status is "new", "open" or blank should remain visible;
# this hides the rows currently marked "Closed" (rows 5, 7, 9, 10) while row 4/6/8/11
# remain visible because their Status still reads "open" at the moment the filter runs.
$lo = $ws.ListObjects.Item(1)
$lo.Range.AutoFilter(1, @("new", "open", ""), 7)

# --- 4. Flip the Status of rows 4 and 8 from "open" to "closed" ---
# This happens after the AutoFilter was applied, so (matching the authored edit)
# these rows stay visible even though they no longer literally match "open".
$ws.Range("A4").Value = "closed"
$ws.Range("A8").Value = "closed"

# --- 5. Fix the typo in the risk description for row 6 (C6): "compated" -> "compared" ---
$ws.Range("C6").Value = "Incomplete and inaccurate data in centralized CCO CSR/Employee database (eWFM/roster raw) - creates inconsistencies in records being created from external sources (quality control system) compared to records created directly in eCL."

# --- 6. Fill in the rest of the new row 12 risk entry ---
$ws.Range("C12").Value = "access to warning information"
$ws.Range("D12").Value = 42060

# --- 7. Fill in new "Date Modified" values for rows 4 and 8 ---
$ws.Range("E4").Value = 42060
$ws.Range("E8").Value = 42060

# --- 8. Update the view: scroll to show row 2 at top, select C13 ---
$excel.ActiveWindow.TopLeftCell = $ws.Range("A2")
$ws.Range("C13").Select()
